# Add the two new footnote entries that were appended to the list in
# column A of Sheet1 (mirrors the two new shared strings / rows in the
# diff: "Task identity" at A65 and "Technique" at A66).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = "Task identity"
$ws.Range("A66").Value = "Technique"

# Reflect the resulting scroll/selection state of the sheet view: the
# window was scrolled down so row 58 is the top visible row, and the
# active/selected cell ended up just past the new last entry (A67).
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A67").Select()
